$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 122, shifting existing rows 122:138 down to 123:139.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new weekly price record.
$ws.Cells.Item(122, 1).Value = 4
$ws.Cells.Item(122, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(122, 3).Value = "Los Lagos"
$ws.Cells.Item(122, 4).Value = 44449
$ws.Cells.Item(122, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(122, 5).Value = 10
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100102
$ws.Cells.Item(122, 8).Value = "Cítricos"
$ws.Cells.Item(122, 9).Value = 100102006
$ws.Cells.Item(122, 10).Value = "Pomelo"
$ws.Cells.Item(122, 11).Value = "Start Ruby"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 160
$ws.Cells.Item(122, 14).Value = 13000
$ws.Cells.Item(122, 15).Value = 13000
$ws.Cells.Item(122, 16).Value = 13000
$ws.Cells.Item(122, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(122, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(122, 19).Value = 929
$ws.Cells.Item(122, 20).Value = 14
